# "saco acentos de los TCs"
# Update the TC (Testigo/NroSiniestro) value in column E row 2,
# and move the selection to E2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the claim number in E2 (keep it as text, preserving the
# existing quote-prefixed text style already applied to the cell).
$ws.Range("E2").Value = "'1120194100405"

# Move / leave the active selection on E2.
$ws.Range("E2").Select() | Out-Null
